# Insert a new data row at row 139 (pushes the existing 139-207 rows down to
# 140-208, matching the growth of the used range from A1:T207 to A1:T208),
# then populate the newly inserted row with its own data. The row that used
# to be at 139 is now at 140 (unchanged), and so on down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(139).Insert()

$ws.Cells.Item(139, 1).Value = 10
$ws.Cells.Item(139, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(139, 3).Value = 'La Araucanía'
$ws.Cells.Item(139, 4).Value = 44460
$ws.Cells.Item(139, 5).Value = 9
$ws.Cells.Item(139, 6).Value = 'Fruta'
$ws.Cells.Item(139, 7).Value = 100108
$ws.Cells.Item(139, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(139, 9).Value = 100108002
$ws.Cells.Item(139, 10).Value = 'Mango'
$ws.Cells.Item(139, 11).Value = 'Sin especificar'
$ws.Cells.Item(139, 12).Value = 'Primera'
$ws.Cells.Item(139, 13).Value = 200
$ws.Cells.Item(139, 14).Value = 9000
$ws.Cells.Item(139, 15).Value = 9000
$ws.Cells.Item(139, 16).Value = 9000
$ws.Cells.Item(139, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(139, 18).Value = 'Brasil'
$ws.Cells.Item(139, 19).Value = 2250
$ws.Cells.Item(139, 20).Value = 4
